# Automatische test-sync: 2025-08-19 21:10:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 25

$ws.Cells.Item($newRow, 1).Value = "Opvolging retour"
$ws.Cells.Item($newRow, 2).Value = "kwaliteit@testbedrijf123.nl"
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 6).Value = "2025-08-19 21:10:46"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges on columns D, G, H, I, J so the
# new row (25) is covered, matching Excel's own "extend formatting" behaviour
# when a new row is appended to a formatted table-like range.
$ranges = @("D2:D24", "G2:G24", "H2:H24", "I2:I24", "J2:J24")
foreach ($rng in $ranges) {
    $col = $rng.Substring(0, 1)
    $newRng = "$col" + "2:$col" + "25"
    $fc = $ws.Range($rng).FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($ws.Range($newRng))
    }
}

# Update the Dashboard summary count for the category that now has one more row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 24
